# One Page Assignment Update
# "turned into PDF form" — consolidate split runs, relocate the stray
# _GoBack bookmark to the trailing empty paragraph, and drop the cached
# PAGE field results from both the header textbox's primary (drawing)
# and fallback (VML) representations.

$d = $word.ActiveDocument

# --- 1. "O" + "verview" -> single run "Overview" (Heading1 paragraph) ---
$p = $d.Paragraphs(2)
$rng = $p.Range
$rng.End = $rng.End - 1
$rng.Delete()
$rng.InsertAfter("Overview")

# --- 2. "Bill the client" + ": A man searching ... suitable." -> single run ---
#     (also removes the stray <w:bookmarkStart/End w:name="_GoBack"/> that
#     sat between the two runs)
$p = $d.Paragraphs(9)
$rng = $p.Range
$rng.End = $rng.End - 1
$sz = $rng.Font.Size
$rng.Delete()
$rng.InsertAfter("Bill the client: A man searching for a capable web developer to create a website to be a boon to his website. He has searched through many profiles to see which ones would be suitable.")
$rng.Font.Size = $sz

# --- 3. Re-add the _GoBack bookmark on the trailing empty paragraph ---
$last = $d.Paragraphs($d.Paragraphs.Count)
$d.Bookmarks.Add("_GoBack", $last.Range)

# --- 4. Strip the cached PAGE field runs from the header's page-number
#     textbox, in both the modern drawing (mc:Choice) and the VML
#     fallback (mc:Fallback) representations. These aren't reachable
#     individually through the Shapes/TextFrame object model (the VML
#     fallback mirrors the drawing and isn't separately editable), so
#     the header part is replaced wholesale via InsertXML with the
#     corrected markup (identical except for the removed field runs).
$sec = $d.Sections(1)
$hdr = $sec.Headers(1)
$hdr.Range.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/header1.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.header+xml"><pkg:xmlData><w:hdr xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:cx="http://schemas.microsoft.com/office/drawing/2014/chartex" xmlns:cx1="http://schemas.microsoft.com/office/drawing/2015/9/8/chartex" xmlns:cx2="http://schemas.microsoft.com/office/drawing/2015/10/21/chartex" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 w16se wp14"><w:p w:rsidR="002B2937" w:rsidRDefault="00956FDF"><w:pPr><w:pStyle w:val="Header"/></w:pPr><w:r><w:rPr><w:noProof/><w:lang w:eastAsia="en-US"/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251659264" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="leftMargin"><wp:align>right</wp:align></wp:positionH><wp:positionV relativeFrom="bottomMargin"><wp:posOffset>0</wp:posOffset></wp:positionV><wp:extent cx="339090" cy="182880"/><wp:effectExtent l="0" t="0" r="3810" b="11430"/><wp:wrapNone/><wp:docPr id="22" name="Text Box 22"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr txBox="1"/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="339090" cy="182880"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln w="6350"><a:noFill/></a:ln><a:effectLst/></wps:spPr><wps:style><a:lnRef idx="0"><a:schemeClr val="accent1"/></a:lnRef><a:fillRef idx="0"><a:schemeClr val="accent1"/></a:fillRef><a:effectRef idx="0"><a:schemeClr val="accent1"/></a:effectRef><a:fontRef idx="minor"><a:schemeClr val="dk1"/></a:fontRef></wps:style><wps:txbx><w:txbxContent><w:p w:rsidR="002B2937" w:rsidRDefault="00956FDF"><w:pPr><w:pStyle w:val="Footer"/></w:pPr></w:p></w:txbxContent></wps:txbx><wps:bodyPr rot="0" spcFirstLastPara="0" vertOverflow="overflow" horzOverflow="overflow" vert="horz" wrap="square" lIns="0" tIns="0" rIns="0" bIns="0" numCol="1" spcCol="0" rtlCol="0" fromWordArt="0" anchor="t" anchorCtr="0" forceAA="0" compatLnSpc="1"><a:prstTxWarp prst="textNoShape"><a:avLst/></a:prstTxWarp><a:spAutoFit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic><wp14:sizeRelH relativeFrom="margin"><wp14:pctWidth>0</wp14:pctWidth></wp14:sizeRelH><wp14:sizeRelV relativeFrom="margin"><wp14:pctHeight>0</wp14:pctHeight></wp14:sizeRelV></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:shapetype id="_x0000_t202" coordsize="21600,21600" o:spt="202" path="m,l,21600r21600,l21600,xe"><v:stroke joinstyle="miter"/><v:path gradientshapeok="t" o:connecttype="rect"/></v:shapetype><v:shape id="Text Box 22" o:spid="_x0000_s1026" type="#_x0000_t202" style="position:absolute;margin-left:-24.5pt;margin-top:0;width:26.7pt;height:14.4pt;z-index:251659264;visibility:visible;mso-wrap-style:square;mso-width-percent:0;mso-height-percent:0;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:right;mso-position-horizontal-relative:left-margin-area;mso-position-vertical:absolute;mso-position-vertical-relative:bottom-margin-area;mso-width-percent:0;mso-height-percent:0;mso-width-relative:margin;mso-height-relative:margin;v-text-anchor:top" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#xA;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#xA;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#xA;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#xA;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#xA;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#xA;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#xA;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#xA;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#xA;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#xA;IQCvuV4gdAIAAFMFAAAOAAAAZHJzL2Uyb0RvYy54bWysVFFP2zAQfp+0/2D5faQtGioVKepATJMQ&#xA;IGDi2XVsGs3xebbbpPv1++wkBbG9MO3Fudx9d7777s5n511j2E75UJMt+fRowpmykqraPpf8++PV&#xA;pzlnIQpbCUNWlXyvAj9ffvxw1rqFmtGGTKU8QxAbFq0r+SZGtyiKIDeqEeGInLIwavKNiPj1z0Xl&#xA;RYvojSlmk8lJ0ZKvnCepQoD2sjfyZY6vtZLxVuugIjMlR24xnz6f63QWyzOxePbCbWo5pCH+IYtG&#xA;1BaXHkJdiijY1td/hGpq6SmQjkeSmoK0rqXKNaCa6eRNNQ8b4VSuBeQEd6Ap/L+w8mZ351ldlXw2&#xA;48yKBj16VF1kX6hjUIGf1oUFYA8OwNhBjz6P+gBlKrvTvklfFMRgB9P7A7spmoTy+Ph0cgqLhGk6&#xA;n83nmf3ixdn5EL8qalgSSu7RvMyp2F2HiEQAHSHpLktXtTG5gcaytuQnx58n2eFggYexCavyKAxh&#xA;UkF94lmKe6MSxth7pUFFzj8p8hCqC+PZTmB8hJTKxlx6jgt0Qmkk8R7HAf+S1Xuc+zrGm8nGg3NT&#xA;W/K5+jdpVz/GlHWPB5Gv6k5i7Nbd0Og1VXv02VO/KcHJqxrduBYh3gmP1UADse7xFoc2BNZpkDjb&#xA;kP/1N33CY2Jh5azFqpU8/NwKrzgz3yxmOe3lKPhRWI+C3TYXBPqneEiczCIcfDSjqD01T3gFVukW&#xA;mISVuKvkcRQvYr/weEWkWq0yCNvnRLy2D06m0KkbabYeuyfh3TCAEZN7Q+MSisWbOeyxeVDcahsx&#xA;jXlIE6E9iwPR2Nw8u8Mrk56G1/8Z9fIWLn8DAAD//wMAUEsDBBQABgAIAAAAIQBckvp92gAAAAMB&#xA;AAAPAAAAZHJzL2Rvd25yZXYueG1sTI/NTsNADITvSLzDykjc6IbyoyjNpkIIeoATKUIc3cTJpmS9&#xA;UXabBp4ewwUulkZjz3zO17Pr1URj6DwbuFwkoIgrX3fcGnjdPl6koEJErrH3TAY+KcC6OD3JMav9&#xA;kV9oKmOrJIRDhgZsjEOmdagsOQwLPxCL1/jRYRQ5troe8SjhrtfLJLnVDjuWBosD3VuqPsqDE4y3&#xA;58Rtvhr77p6wCaXdTpuHvTHnZ/PdClSkOf4tww++3EAhTDt/4Dqo3oA8En+neDdX16B2BpZpCrrI&#xA;9X/24hsAAP//AwBQSwECLQAUAAYACAAAACEAtoM4kv4AAADhAQAAEwAAAAAAAAAAAAAAAAAAAAAA&#xA;W0NvbnRlbnRfVHlwZXNdLnhtbFBLAQItABQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAAAAAAAAA&#xA;AAAAAC8BAABfcmVscy8ucmVsc1BLAQItABQABgAIAAAAIQCvuV4gdAIAAFMFAAAOAAAAAAAAAAAA&#xA;AAAAAC4CAABkcnMvZTJvRG9jLnhtbFBLAQItABQABgAIAAAAIQBckvp92gAAAAMBAAAPAAAAAAAA&#xA;AAAAAAAAAM4EAABkcnMvZG93bnJldi54bWxQSwUGAAAAAAQABADzAAAA1QUAAAAA&#xA;" filled="f" stroked="f" strokeweight=".5pt"><v:textbox style="mso-fit-shape-to-text:t" inset="0,0,0,0"><w:txbxContent><w:p w:rsidR="002B2937" w:rsidRDefault="00956FDF"><w:pPr><w:pStyle w:val="Footer"/></w:pPr></w:p></w:txbxContent></v:textbox><w10:wrap anchorx="margin" anchory="margin"/></v:shape></w:pict></mc:Fallback></mc:AlternateContent></w:r></w:p></w:hdr></pkg:xmlData></pkg:part></pkg:package>')

Write-Host "Edits applied"
